# The "personal_info" sheet's A1, B1 and C1 cells had their inline-string
# text content deleted: each cell stays present on the row (still a text
# cell) but now holds an empty string, while D1 ("1") is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("personal_info")

$target = $ws.Range("A1:C1")

# Assigning a bare "'" forces each cell to stay text-typed with empty
# content (rather than reverting to a blank/numeric cell, which Excel
# would otherwise drop entirely on save).
$target.Value = "'"

# Re-apply the default style so no stray quote-prefix formatting is left
# behind on the cells themselves.
$target.Style = "Normal"
